# Atualizar fontes de dados
# Drop the oldest month (old row 2) so every subsequent row shifts up by one
# position: year/month move forward by a month, and each row's maturation
# label is recomputed relative to the new starting month. The now-stale
# numeric measures (n_hours .. profitLoss) are reset to 0 until the new
# source data is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old 2025-08 / m0); everything below shifts up,
# which naturally advances year/month for every remaining row by one month.
$ws.Rows.Item(2).Delete()

# Recompute the "maturation" label per row: the first four remaining rows
# are the new near-term buckets m0..m3, after which rows are grouped into
# annual buckets a0, a1, a2, ... by calendar year.
$labels = @("m0", "m1", "m2", "m3")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value2 = $labels[$i]
}

$lastRow = $ws.UsedRange.Rows.Count
$yearOffset = -1
$prevYear = -1
for ($r = 6; $r -le $lastRow; $r++) {
    $yr = $ws.Cells.Item($r, 1).Value2
    if ($yr -ne $prevYear) {
        $yearOffset = $yearOffset + 1
        $prevYear = $yr
    }
    $ws.Cells.Item($r, 3).Value2 = "a" + $yearOffset
}

# The measures (n_hours .. profitLoss) are stale for the refreshed window
# until new source data lands; zero them out for every remaining data row.
$numRows = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range($ws.Cells.Item(2, 4), $ws.Cells.Item($numRows, 17))
$dataRange.Value2 = 0
